{"js": "// worked on public_benefits overflow\n//\n// The addendum's caption line used to read:\n//     {{ users[0] }} v. {{ other_parties[0] }}\n// It is collapsed down to a single merge field:\n//     {{ case_name }}\n//\n// We locate that text, then replace the matched range's contents with the\n// new run layout via insertOoxml. Because the match only spans the\n// paragraph's content (not its trailing paragraph mark), the host\n// paragraph's own properties (style, paraId, rsids, ...) are left\n// completely untouched - only the runs inside it change, exactly as the\n// diff shows.\n\nconst results = context.document.body.search(\n  \"{{ users[0] }} v. {{ other_parties[0] }}\",\n  { matchCase: true, matchWildcards: false }\n);\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n\n  const newOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>{</w:t></w:r>\n            <w:r w:rsidRPr=\"00F54495\"><w:t>{</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r w:rsidRPr=\"00F54495\"><w:t>case_name</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> }</w:t></w:r>\n            <w:r w:rsidRPr=\"00F54495\"><w:t>}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n  target.insertOoxml(newOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# worked on public_benefits overflow\n#\n# The addendum's caption line used to read:\n#     {{ users[0] }} v. {{ other_parties[0] }}\n# It is collapsed down to a single merge field:\n#     {{ case_name }}\n#\n# We locate the paragraph by its distinctive text, then replace only the\n# *content* of that paragraph (i.e. everything except the trailing\n# paragraph-mark) with the new run layout. Leaving the paragraph mark out\n# of the replaced range means the paragraph's own properties (style,\n# paraId, rsids, ...) are left completely untouched - only the runs inside\n# it change, exactly as the diff shows.\n\n$d = $word.ActiveDocument\n\n# NOTE: Find.Execute() mutates the Range it was obtained from in place -\n# it does NOT return a new Range. So we must keep a handle to that same\n# Range object ($target) rather than re-querying $d.Content afterwards.\n$target = $d.Content\n$find = $target.Find\n$find.ClearFormatting()\n$find.Text = \"{{ users[0] }} v. {{ other_parties[0] }}\"\n$find.MatchWildcards = $false\n\nif ($find.Execute()) {\n    # $target now spans exactly the matched text (Find does not include\n    # the trailing paragraph-mark), so re-wrapping it in a fresh Range\n    # keeps the host paragraph (pPr, paraId, rsids, ...) untouched - only\n    # the runs inside it are replaced.\n    $contentRange = $d.Range($target.Start, $target.End)\n\n    $newXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>{</w:t></w:r>\n            <w:r w:rsidRPr=\"00F54495\"><w:t>{</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r w:rsidRPr=\"00F54495\"><w:t>case_name</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> }</w:t></w:r>\n            <w:r w:rsidRPr=\"00F54495\"><w:t>}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n    $contentRange.InsertXML($newXml)\n}\n"}
